# Updates cryptocurrency price/volume data per upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.644.03"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.170.09"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Formula = "'226.63"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Formula = "'0.626"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").Formula = "'63.20"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Formula = "'0.391"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Formula = "'0.0851"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Formula = "'15.92"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").Value = "2.489.58"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Formula = "'21.80"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").Formula = "'0.811"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Formula = "'5.50"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").Value = "2.167.95"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "39.587.15"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "0.0₃0920"
$ws.Range("E19").Value = "  +7.62%  "
$ws.Range("D20").Formula = "'71.72"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("D22").Formula = "'229.60"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Formula = "'2.32"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("D25").Formula = "'2.34"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").Formula = "'171.00"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").Formula = "'9.49"
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Formula = "'1.45"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Formula = "'19.82"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("E31").Value = "  +3.66%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Formula = "'4.52"
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("D34").Formula = "'4.69"
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("D35").Formula = "'6.97"
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("D36").Formula = "'0.0617"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Formula = "'3.82"
$ws.Range("E37").Value = "  +6.51%  "
$ws.Range("D38").Formula = "'2.40"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Formula = "'4.91"
$ws.Range("E40").Value = "  +17.75%  "
# Row 41 now lists VeChain (previously Aave)
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Formula = "'0.0230"
$ws.Range("E41").Value = "  -1.26%  "
# Row 42 now lists Aave (previously VeChain)
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Formula = "'102.76"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Formula = "'17.70"
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("D44").Value = "1.512.83"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").Formula = "'0.0922"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").Formula = "'0.000196"
$ws.Range("E50").Value = "  +32.79%  "
# Row 51 now lists MultiversX (previously RocketPoolETH)
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Formula = "'49.60"
$ws.Range("E51").Value = "  +6.46%  "
